$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix bad data in column B (age bracket: "y"/"o") ---
$ws.Range("B12").Value = "y"
$ws.Range("B13").Value = "o"
$ws.Range("B15").Value = "y"
$ws.Range("B16").Value = "o"

# --- Highlight the corrected column plus a new helper column with green fill ---
$green = 5287936  # RGB(0, 176, 80) packed as BGR OLE color
$ws.Range("B2:B15").Interior.Color = $green
$ws.Range("D2:D15").Interior.Color = $green

# --- Match the new selection left behind on the sheet ---
$ws.Range("C1:D17").Select()
